$wb = $excel.ActiveWorkbook

# The new "Croatia" test-data sheet is a copy of the existing "Spain" sheet
# (same layout/styles), placed as the last tab and becoming the active tab.
$spain = $wb.Worksheets.Item("Spain")
$spain.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))

$croatia = $wb.Worksheets.Item($wb.Worksheets.Count)
$croatia.Name = "Croatia"

# Update the two market-specific cells (order matters for shared-string index
# allocation: NGC ticket reference first, then the market name).
$croatia.Range("B4").Value = "NGC-3139/T2485"
$croatia.Range("B2").Value = "Croatia Market"

$croatia.Select()
